$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# Update the time_taken column (F) on the data sheet with refreshed timestamps
$data.Range("F2").Value = "2021-10-05 14:21:26.742871"
$data.Range("F3").Value = "2021-10-05 14:21:26.742879"
$data.Range("F4").Value = "2021-10-05 14:21:26.742882"
$data.Range("F5").Value = "2021-10-05 14:21:26.742885"
$data.Range("F6").Value = "2021-10-05 14:21:26.742887"
$data.Range("F7").Value = "2021-10-05 14:21:26.742890"
$data.Range("F8").Value = "2021-10-05 14:21:26.742892"
$data.Range("F9").Value = "2021-10-05 14:21:26.742895"
$data.Range("F10").Value = "2021-10-05 14:21:26.742898"
$data.Range("F11").Value = "2021-10-05 14:21:26.742900"
$data.Range("F12").Value = "2021-10-05 14:21:26.742903"
$data.Range("F13").Value = "2021-10-05 14:21:26.742905"
$data.Range("F14").Value = "2021-10-05 14:21:26.742907"
$data.Range("F15").Value = "2021-10-05 14:21:26.742910"
$data.Range("F16").Value = "2021-10-05 14:21:26.742913"
$data.Range("F17").Value = "2021-10-05 14:21:26.742915"
$data.Range("F18").Value = "2021-10-05 14:21:26.742918"
$data.Range("F19").Value = "2021-10-05 14:21:26.742920"
$data.Range("F20").Value = "2021-10-05 14:21:26.742923"
$data.Range("F21").Value = "2021-10-05 14:21:26.742925"
$data.Range("F22").Value = "2021-10-05 14:21:26.742928"
$data.Range("F23").Value = "2021-10-05 14:21:26.742930"

# Add a new "metadata" worksheet right after the "data" sheet
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Copy the header / index-column formatting from the "data" sheet so the
# new sheet reuses the same bold/bordered/centered style
$data.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Long QT syndrome"
$meta.Range("C2").Value = 76
# data_version is stored as text ("2.23"), not a number, so force text format first
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "2.23"
$meta.Range("E2").Value = "2021-09-28T09:41:48.458571Z"
$meta.Range("F2").Value = "2021-10-05 14:21:26.739164"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/76/?format=json"
